$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actors")

# Update row 2 ("Cube") values: LocationZ, RotationZ, Scale
$ws.Range("B2").Value = 200
$ws.Range("C2").Value = 30
$ws.Range("D2").Value = 2

# Update the selected/active cell on the sheet to match the saved view state
$ws.Activate()
$ws.Range("C7").Select()
